# Reproduce the "modified code in excel reader" edit:
#  - Sheet2!B2 used to hold the formula =Sheet1!B3*2; it now holds a
#    plain text value of a single space character.
#  - The cursor/selection left behind in each sheet moved (Sheet1 -> A3,
#    Sheet2 -> C6), and Sheet2 is the active (visible) sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1: leave the selection on A3 (was B2).
$ws1.Activate()
$ws1.Range("A3").Select()

# Sheet2: replace the formula in B2 with a literal single space, then
# leave the selection on C6 (was B6) with Sheet2 as the active sheet.
$ws2.Activate()
$ws2.Range("B2").Value = " "
$ws2.Range("C6").Select()
